# Logged Week 16 and performed season sim from Week 17
# Adds two new players (J.Moore to WR, M.Ffrench to RB) as new rows of
# all-zero stats on the RB and WR sheets.

$wb = $excel.ActiveWorkbook

# --- WR sheet: append J.Moore as a new row (row 8) ---
# (added first so its shared-string entry precedes M.Ffrench's, matching
# the new-player logging order for the week)
$wr = $wb.Worksheets.Item("WR")
$wrRow = 8
$wr.Cells.Item($wrRow, 1).Value = "J.Moore"
for ($col = 2; $col -le 10; $col++) {
    $wr.Cells.Item($wrRow, $col).Value = 0
}
$wr.Range("J9").Select()

# --- RB sheet: append M.Ffrench as a new row (row 7) ---
$rb = $wb.Worksheets.Item("RB")
$rbRow = 7
$rb.Cells.Item($rbRow, 1).Value = "M.Ffrench"
for ($col = 2; $col -le 10; $col++) {
    $rb.Cells.Item($rbRow, $col).Value = 0
}
$rb.Range("J8").Select()
